# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 349, pushing the existing
# rows 349-380 down to 350-381 (dimension grows from A1:R380 to A1:R381).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the existing rows down by inserting a fresh row 349.
$ws.Rows.Item(349).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(349, 1).Value  = 3
$ws.Cells.Item(349, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(349, 3).Value  = "Coquimbo"
$ws.Cells.Item(349, 4).Value  = 44769
$ws.Cells.Item(349, 5).Value  = 5
$ws.Cells.Item(349, 6).Value  = 100112009
$ws.Cells.Item(349, 7).Value  = "Acelga"
$ws.Cells.Item(349, 8).Value  = "Sin especificar"
$ws.Cells.Item(349, 9).Value  = "Primera"
$ws.Cells.Item(349, 10).Value = 260
$ws.Cells.Item(349, 11).Value = 3000
$ws.Cells.Item(349, 12).Value = 3300
$ws.Cells.Item(349, 13).Value = 3162
$ws.Cells.Item(349, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(349, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(349, 16).Value = 527
$ws.Cells.Item(349, 17).Value = 6
$ws.Cells.Item(349, 18).Value = "Hortaliza"
